$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 184, shifting the existing rows 184-187 down to 185-188.
$ws.Rows.Item(184).Insert()

# Fill the new row 184 with the new weekly record.
$ws.Cells.Item(184, 1).Value = 3
$ws.Cells.Item(184, 2).Value = "Femacal de La Calera"
$ws.Cells.Item(184, 3).Value = "Coquimbo"
$ws.Cells.Item(184, 4).Value = 44448
$ws.Cells.Item(184, 5).Value = 5
$ws.Cells.Item(184, 6).Value = 100114013
$ws.Cells.Item(184, 7).Value = "Zanahoria"
$ws.Cells.Item(184, 8).Value = "Sin especificar"
$ws.Cells.Item(184, 9).Value = "Primera"
$ws.Cells.Item(184, 10).Value = 500
$ws.Cells.Item(184, 11).Value = 5000
$ws.Cells.Item(184, 12).Value = 5500
$ws.Cells.Item(184, 13).Value = 5260
$ws.Cells.Item(184, 14).Value = "$/saco 20 kilos"
$ws.Cells.Item(184, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(184, 16).Value = 263
$ws.Cells.Item(184, 17).Value = 20
$ws.Cells.Item(184, 18).Value = "Hortaliza"
